$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the D/F values between rows 15 and 16
$d15 = $ws.Range("D15").Value2
$f15 = $ws.Range("F15").Value2
$d16 = $ws.Range("D16").Value2
$f16 = $ws.Range("F16").Value2

$ws.Range("D15").Value = $d16
$ws.Range("F15").Value = $f16
$ws.Range("D16").Value = $d15
$ws.Range("F16").Value = $f15
